$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("optimal models")
$ws2 = $wb.Worksheets.Item("relative importance")

# --- Sheet1 ("optimal models"): add row 16 for LESC (longfin escolar) ---
$ws1.Range("A16").Value = "LESC"
$ws1.Range("B16").Value = 1
$ws1.Range("C16").Value = 0.75
$ws1.Range("D16").Value = 0.005
$ws1.Range("E16").Value = 6450
$ws1.Range("F16").Value = 0.8857
$ws1.Range("G16").Value = 0.8011
$ws1.Range("H16").Value = 0.8847
$ws1.Range("I16").Value = 0.127967

$ws1.Range("I16").Select()

# --- Sheet2 ("relative importance"): add row 16 for LESC ---
$ws2.Range("A16").Value = "LESC"
$ws2.Range("B16").Value = 3.2421185
$ws2.Range("C16").Value = 8.4286494
$ws2.Range("D16").Value = 1.5517163
$ws2.Range("E16").Value = 2.3250592
$ws2.Range("F16").Value = 2.7232842
$ws2.Range("G16").Value = 6.7300705
$ws2.Range("H16").Value = 16.5446309
$ws2.Range("I16").Value = 4.2541132
$ws2.Range("J16").Value = 6.3973210
$ws2.Range("K16").Value = 0.1912612
$ws2.Range("L16").Value = 1.4112623
$ws2.Range("M16").Value = 6.3112581
$ws2.Range("N16").Value = 4.9988838
$ws2.Range("O16").Value = 4.9516112
$ws2.Range("P16").Value = 19.4626212
$ws2.Range("Q16").Value = 9.1104019
$ws2.Range("R16").Value = 1.3657369

$ws2.Range("K16").Select()

# Restore "optimal models" as the active/selected sheet (as in the source workbook)
$ws1.Activate()
$ws1.Range("I16").Select()
